$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.352.45'
$ws.Range('E2').Value = '  +4.24%  '

$ws.Range('D3').Value = '2.434.39'
$ws.Range('E3').Value = '  +3.18%  '

$ws.Range('E4').Value = '  +0.00%  '

$c = $ws.Range('D5')
$c.Value = "'556.49"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.45%  '

$c = $ws.Range('D6')
$c.Value = "'139.54"
$c.Style = 'Normal'

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  +3.01%  '

$ws.Range('E9').Value = '  +4.99%  '

$ws.Range('E11').Value = '  +1.37%  '

$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('E13').Value = '  +5.51%  '

$ws.Range('D14').Value = '2.866.18'
$ws.Range('E14').Value = '  +3.13%  '

$ws.Range('D15').Value = '60.270.62'
$ws.Range('E15').Value = '  +4.16%  '

$ws.Range('E16').Value = '  +4.29%  '

$ws.Range('D17').Value = '2.432.70'
$ws.Range('E17').Value = '  +3.45%  '

$ws.Range('E18').Value = '  +5.87%  '

$ws.Range('E19').Value = '  +3.04%  '

$c = $ws.Range('D20')
$c.Value = "'334.11"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.97%  '

$ws.Range('E21').Value = '  +0.93%  '

$c = $ws.Range('D22')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '

$c = $ws.Range('D23')
$c.Value = "'65.35"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +4.30%  '

$ws.Range('E24').Value = '  +3.45%  '

$ws.Range('E25').Value = '  +1.92%  '

$ws.Range('E26').Value = '  +0.01%  '

$c = $ws.Range('D27')
$c.Value = "'1.35"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.17%  '

$ws.Range('E28').Value = '  +6.93%  '

$ws.Range('E29').Value = '  +2.15%  '

$c = $ws.Range('D30')
$c.Value = "'6.35"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +3.29%  '

$c = $ws.Range('D31')
$c.Value = "'169.29"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.59%  '

$c = $ws.Range('D32')
$c.Value = "'1.05"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +1.05%  '

$ws.Range('E33').Value = '  +1.95%  '

$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('E35').Value = '  +6.24%  '

$ws.Range('E36').Value = '  +0.33%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('E38').Value = '  +0.34%  '

$ws.Range('E39').Value = '  +11.40%  '

$c = $ws.Range('D40')
$c.Value = "'39.87"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.18%  '

$c = $ws.Range('D41')
$c.Value = "'322.34"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +11.64%  '

$ws.Range('E42').Value = '  +1.58%  '

$c = $ws.Range('D43')
$c.Value = "'141.31"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.11%  '

$ws.Range('E44').Value = '  +3.69%  '

$c = $ws.Range('D45')
$c.Value = "'0.0962"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.68%  '

$ws.Range('B46').Value = 'Polygon'
$ws.Range('C46').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D46')
$c.Value = "'0.419"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +9.21%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D47')
$c.Value = "'19.62"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.83%  '

$ws.Range('E48').Value = '  +1.46%  '

$ws.Range('E49').Value = '  +2.27%  '

$ws.Range('E50').Value = '  +2.76%  '
